$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, copying the same formatting used by the
# other header cells (bold, centered, bordered) from G1 ("sum").
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# New "Save" value for the single data row.
$ws.Range("H2").Value = 1
